$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# hunk 0: row 129
$ws.Range("H129").Value = 257962.33
$ws.Range("I129").Value = 379.88
$ws.Range("J129").Value = 717931
$ws.Range("K129").Value = 1139.64
$ws.Range("L129").Value = 2153793
$ws.Range("M129").Value = 3860.36
$ws.Range("N129").Value = -2163793
# hunk 1: row 137
$ws.Range("H137").Value = 10785272
$ws.Range("I137").Value = 962.55316
$ws.Range("J137").Value = 56863690
$ws.Range("K137").Value = 2887.65948
$ws.Range("L137").Value = 170591070
$ws.Range("M137").Value = -337.6594800000003
$ws.Range("N137").Value = -170596170

$ws = $wb.Worksheets.Item("ARM")
# hunk 2: row 32
$ws.Range("H32").Value = 7156383.5
$ws.Range("I32").Value = 13684.75
$ws.Range("J32").Value = 22740454
$ws.Range("K32").Value = 13684.75
$ws.Range("L32").Value = 22740454
$ws.Range("M32").Value = -13397.75
$ws.Range("N32").Value = -22741028
# hunk 3: row 74
$ws.Range("H74").Value = 30668324
$ws.Range("I74").Value = 27028094
$ws.Range("J74").Value = 41028980
$ws.Range("K74").Value = 27028094
$ws.Range("L74").Value = 41028980
$ws.Range("M74").Value = -27027220
$ws.Range("N74").Value = -41030728
# hunk 4: row 77
$ws.Range("H77").Value = 30668324
$ws.Range("I77").Value = 27028094
$ws.Range("J77").Value = 41028980
$ws.Range("K77").Value = 135140470
$ws.Range("L77").Value = 205144900
$ws.Range("M77").Value = -135136102
$ws.Range("N77").Value = -205153636
# hunk 5: row 132
$ws.Range("H132").Value = 13429294
$ws.Range("I132").Value = 17395496
$ws.Range("J132").Value = 4276520.5
$ws.Range("K132").Value = 52186488
$ws.Range("L132").Value = 12829561.5
$ws.Range("M132").Value = -52183958
$ws.Range("N132").Value = -12834621.5

$ws = $wb.Worksheets.Item("BSM")
# hunk 6: row 134
$ws.Range("H134").Value = 13737396
$ws.Range("I134").Value = 15152462
$ws.Range("K134").Value = 45457386
$ws.Range("M134").Value = -45454851

$ws = $wb.Worksheets.Item("CRP")
# hunk 7: row 31
$ws.Range("H31").Value = 1693625
$ws.Range("I31").Value = 1900.4348
$ws.Range("J31").Value = 4472886.5
$ws.Range("K31").Value = 1900.4348
$ws.Range("L31").Value = 4472886.5
$ws.Range("M31").Value = -1605.4348
$ws.Range("N31").Value = -4473476.5
# hunk 8: row 34
$ws.Range("H34").Value = 1693625
$ws.Range("I34").Value = 1900.4348
$ws.Range("J34").Value = 4472886.5
$ws.Range("K34").Value = 1900.4348
$ws.Range("L34").Value = 4472886.5
$ws.Range("M34").Value = -1698.4348
$ws.Range("N34").Value = -4473290.5
# hunk 9: row 51
$ws.Range("H51").Value = 9400.666999999999
$ws.Range("J51").Value = 9400.666999999999
$ws.Range("L51").Value = 9400.666999999999
$ws.Range("N51").Value = -10872.667
# hunk 10: row 58
$ws.Range("H58").Value = 970024.8
$ws.Range("I58").Value = 3737.6875
$ws.Range("J58").Value = 3031437.5
$ws.Range("K58").Value = 3737.6875
$ws.Range("L58").Value = 3031437.5
$ws.Range("M58").Value = -3534.6875
$ws.Range("N58").Value = -3031843.5
# hunk 11: row 59
$ws.Range("H59").Value = 16069.5
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 16069.5
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 16069.5
$ws.Range("N59").Value = -18359.5
$ws.Range("M59").ClearContents()
# hunk 12: row 60
$ws.Range("H60").Value = 7093.6665
$ws.Range("J60").Value = 8200.5
$ws.Range("L60").Value = 8200.5
$ws.Range("N60").Value = -9222.5
# hunk 13: row 61
$ws.Range("H61").Value = 9400.666999999999
$ws.Range("J61").Value = 9400.666999999999
$ws.Range("L61").Value = 9400.666999999999
$ws.Range("N61").Value = -10096.667
# hunk 14: row 68
$ws.Range("H68").Value = 17685.428
$ws.Range("J68").Value = 19299.666
$ws.Range("L68").Value = 19299.666
$ws.Range("N68").Value = -20797.666
# hunk 15: row 71
$ws.Range("H71").Value = 17685.428
$ws.Range("J71").Value = 19299.666
$ws.Range("L71").Value = 57898.99800000001
$ws.Range("N71").Value = -65386.99800000001
# hunk 16: row 74
$ws.Range("H74").Value = 17305.572
$ws.Range("J74").Value = 18556.5
$ws.Range("L74").Value = 18556.5
$ws.Range("N74").Value = -20304.5
# hunk 17: row 77
$ws.Range("H77").Value = 17305.572
$ws.Range("J77").Value = 18556.5
$ws.Range("L77").Value = 55669.5
$ws.Range("N77").Value = -64405.5
# hunk 18: row 107
$ws.Range("H107").Value = 929.6875
$ws.Range("I107").Value = 284.3
$ws.Range("J107").Value = 1223.0454
$ws.Range("K107").Value = 284.3
$ws.Range("L107").Value = 1223.0454
$ws.Range("M107").Value = 1635.7
$ws.Range("N107").Value = -5063.0454
# hunk 19: row 132
$ws.Range("H132").Value = 1537.4314
$ws.Range("I132").Value = 1202.8684
$ws.Range("J132").Value = 2515.3845
$ws.Range("K132").Value = 3608.6052
$ws.Range("L132").Value = 7546.1535
$ws.Range("M132").Value = -1078.6052
$ws.Range("N132").Value = -12606.1535
# hunk 20: row 134
$ws.Range("H134").Value = 934837.9399999999
$ws.Range("I134").Value = 5252.552
$ws.Range("J134").Value = 2860407.8
$ws.Range("K134").Value = 15757.656
$ws.Range("L134").Value = 8581223.399999999
$ws.Range("M134").Value = -13222.656
$ws.Range("N134").Value = -8586293.399999999
# hunk 21: row 136
$ws.Range("H136").Value = 970024.8
$ws.Range("I136").Value = 3737.6875
$ws.Range("J136").Value = 3031437.5
$ws.Range("K136").Value = 11213.0625
$ws.Range("L136").Value = 9094312.5
$ws.Range("M136").Value = -8663.0625
$ws.Range("N136").Value = -9099412.5
# hunk 22: row 140
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
# hunk 23: row 141
$ws.Range("H141").Value = 392393.6
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 392393.6
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 392393.6
$ws.Range("N141").Value = -402753.6
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# hunk 24: row 132
$ws.Range("H132").Value = 5643470.5
$ws.Range("I132").Value = 5503738.5
$ws.Range("J132").Value = 6062665.5
$ws.Range("K132").Value = 16511215.5
$ws.Range("L132").Value = 18187996.5
$ws.Range("M132").Value = -16508685.5
$ws.Range("N132").Value = -18193056.5

$ws = $wb.Worksheets.Item("LTW")
# hunk 25: row 132
$ws.Range("H132").Value = 1907389.1
$ws.Range("I132").Value = 2344417.2
$ws.Range("J132").Value = 3195.0715
$ws.Range("K132").Value = 7033251.600000001
$ws.Range("L132").Value = 9585.2145
$ws.Range("M132").Value = -7030721.600000001
$ws.Range("N132").Value = -14645.2145

$ws = $wb.Worksheets.Item("WVR")
# hunk 26: row 5
$ws.Range("H5").Value = 3650
$ws.Range("J5").Value = 3650
$ws.Range("L5").Value = 3650
$ws.Range("N5").Value = -3874
# hunk 27: row 132
$ws.Range("H132").Value = 637084.8
$ws.Range("I132").Value = 2004
$ws.Range("K132").Value = 6012
$ws.Range("M132").Value = -3482
# hunk 28: row 136
$ws.Range("H136").Value = 3133.2056
$ws.Range("I136").Value = 2093.0571
$ws.Range("J136").Value = 4091.2368
$ws.Range("K136").Value = 6279.1713
$ws.Range("L136").Value = 12273.7104
$ws.Range("M136").Value = -3729.1713
$ws.Range("N136").Value = -17373.7104
